$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.03
$ws.Range("O2").Value = 1.19
$ws.Range("S2").Value = 2.7
$ws.Range("T2").Value = 1.41

# Row 3
$ws.Range("M3").Value = 1.08
$ws.Range("O3").Value = 1.46
$ws.Range("P3").Value = 2.57
$ws.Range("T3").Value = 1.13

# Row 4
$ws.Range("M4").Value = 1.1
$ws.Range("O4").Value = 1.54
$ws.Range("T4").Value = 1.1

# Row 5
$ws.Range("M5").Value = 1.07
$ws.Range("O5").Value = 1.47
$ws.Range("T5").Value = 1.13
$ws.Range("AP5").Value = 1.93
$ws.Range("AQ5").Value = 1.93

# Row 6
$ws.Range("M6").Value = 1.1
$ws.Range("O6").Value = 1.54
$ws.Range("T6").Value = 1.1

# Row 7
$ws.Range("G7").Value = 2.25
$ws.Range("H7").Value = 2.75
$ws.Range("I7").Value = 4
$ws.Range("J7").Value = 3.1
$ws.Range("L7").Value = 4.75
$ws.Range("M7").Value = 1.14
$ws.Range("N7").Value = 5.5
$ws.Range("Q7").Value = 3.1
$ws.Range("R7").Value = 1.36
$ws.Range("S7").Value = 6.5
$ws.Range("T7").Value = 1.11
$ws.Range("U7").Value = 1.67
$ws.Range("V7").Value = 2.1
$ws.Range("Y7").Value = 5
$ws.Range("Z7").Value = 8.5
$ws.Range("AL7").Value = 17

# Row 9
$ws.Range("Q9").Value = 2.15
$ws.Range("R9").Value = 1.67
$ws.Range("S9").Value = 3.75
$ws.Range("T9").Value = 1.25

# Row 13
$ws.Range("G13").Value = 2.6
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 2.57
$ws.Range("J13").Value = 3.15
$ws.Range("L13").Value = 3.1
$ws.Range("N13").Value = 7.3
$ws.Range("Q13").Value = 2.05
$ws.Range("R13").Value = 1.72
$ws.Range("U13").Value = 1.42
$ws.Range("V13").Value = 2.75
$ws.Range("Y13").Value = 7.7
$ws.Range("AB13").Value = 30
$ws.Range("AE13").Value = 7.3
$ws.Range("AF13").Value = 6.8
$ws.Range("AJ13").Value = 7.7
$ws.Range("AL13").Value = 10.5
$ws.Range("AN13").Value = 24

# Row 14
$ws.Range("G14").Value = 4.35
$ws.Range("I14").Value = 1.8
$ws.Range("J14").Value = 4.7
$ws.Range("K14").Value = 2.1
$ws.Range("L14").Value = 2.37
$ws.Range("S14").Value = 3.8
$ws.Range("U14").Value = 1.45
$ws.Range("Z14").Value = 24
$ws.Range("AD14").Value = 65
$ws.Range("AM14").Value = 15

# Row 18
$ws.Range("G18").Value = 1.57
$ws.Range("H18").Value = 3.7
$ws.Range("I18").Value = 5.4
$ws.Range("J18").Value = 2.12
$ws.Range("K18").Value = 2.15
$ws.Range("L18").Value = 5.4
$ws.Range("Q18").Value = 1.75
$ws.Range("R18").Value = 1.85
$ws.Range("S18").Value = 2.77
$ws.Range("T18").Value = 1.34
$ws.Range("W18").Value = 1.78
$ws.Range("X18").Value = 1.82
$ws.Range("Y18").Value = 6.9
$ws.Range("Z18").Value = 7.5
$ws.Range("AB18").Value = 11.75
$ws.Range("AC18").Value = 12.5
$ws.Range("AD18").Value = 25
$ws.Range("AF18").Value = 7.3
$ws.Range("AG18").Value = 16.5
$ws.Range("AH18").Value = 80
$ws.Range("AI18").Value = 600
$ws.Range("AJ18").Value = 14
$ws.Range("AK18").Value = 32
$ws.Range("AL18").Value = 17
$ws.Range("AM18").Value = 110
$ws.Range("AN18").Value = 60
$ws.Range("AO18").Value = 60

# Row 20
$ws.Range("G20").Value = 1.87
$ws.Range("H20").Value = 3.3
$ws.Range("I20").Value = 4.1
$ws.Range("N20").Value = 8.5
$ws.Range("Z20").Value = 8.5

# Row 22
$ws.Range("G22").Value = 2.65
$ws.Range("I22").Value = 2.35

# Row 23
$ws.Range("I23").Value = 1.81

# Row 24
$ws.Range("I24").Value = 1.71
$ws.Range("O24").Value = 1.11
$ws.Range("P24").Value = 6.5
$ws.Range("S24").Value = 1.91
$ws.Range("T24").Value = 1.8

# Row 25
$ws.Range("G25").Value = 2.37
$ws.Range("N25").Value = 9

# Row 26
$ws.Range("G26").Value = 1.96

# Row 27
$ws.Range("I27").Value = 2.87

# Row 29
$ws.Range("G29").Value = 1.63

# Row 30
$ws.Range("G30").Value = 2.38
$ws.Range("H30").Value = 3.3
$ws.Range("I30").Value = 2.8
$ws.Range("J30").Value = 2.88
$ws.Range("K30").Value = 2.2
$ws.Range("L30").Value = 3.25
$ws.Range("M30").Value = 1.01
$ws.Range("N30").Value = 11
$ws.Range("O30").Value = 1.22
$ws.Range("P30").Value = 3.75
$ws.Range("Q30").Value = 1.83
$ws.Range("R30").Value = 1.98
$ws.Range("S30").Value = 3
$ws.Range("U30").Value = 1.36
$ws.Range("V30").Value = 3
$ws.Range("W30").Value = 1.67
$ws.Range("X30").Value = 2.1
$ws.Range("Y30").Value = 9.5
$ws.Range("Z30").Value = 12
$ws.Range("AB30").Value = 23
$ws.Range("AC30").Value = 19
$ws.Range("AD30").Value = 26
$ws.Range("AE30").Value = 11
$ws.Range("AF30").Value = 6.5
$ws.Range("AG30").Value = 13
$ws.Range("AH30").Value = 41
$ws.Range("AI30").Value = 151
$ws.Range("AJ30").Value = 10
$ws.Range("AL30").Value = 11
$ws.Range("AM30").Value = 29
$ws.Range("AN30").Value = 21
$ws.Range("AO30").Value = 29

# Row 31
$ws.Range("G31").Value = 1.73
$ws.Range("I31").Value = 4.75
$ws.Range("K31").Value = 2.2
$ws.Range("L31").Value = 4.75
$ws.Range("M31").Value = 1.02
$ws.Range("N31").Value = 10
$ws.Range("O31").Value = 1.22
$ws.Range("P31").Value = 3.75
$ws.Range("Q31").Value = 1.9
$ws.Range("R31").Value = 1.9
$ws.Range("S31").Value = 3.25
$ws.Range("T31").Value = 1.3
$ws.Range("U31").Value = 1.4
$ws.Range("V31").Value = 2.75
$ws.Range("W31").Value = 1.83
$ws.Range("X31").Value = 1.83
$ws.Range("Y31").Value = 7.5
$ws.Range("Z31").Value = 8
$ws.Range("AA31").Value = 9
$ws.Range("AB31").Value = 13
$ws.Range("AC31").Value = 15
$ws.Range("AE31").Value = 10
$ws.Range("AF31").Value = 7
$ws.Range("AH31").Value = 51
$ws.Range("AK31").Value = 26
$ws.Range("AM31").Value = 51
$ws.Range("AN31").Value = 41
$ws.Range("AO31").Value = 41

# Row 33
$ws.Range("G33").Value = 2.2
$ws.Range("I33").Value = 3
$ws.Range("J33").Value = 2.63
$ws.Range("L33").Value = 3.25
$ws.Range("M33").Value = 1.01
$ws.Range("N33").Value = 19
$ws.Range("O33").Value = 1.1
$ws.Range("T33").Value = 1.69
$ws.Range("AA33").Value = 9.5
